$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 0.1839054330529666
$ws.Range("C2").Value = 0.2158649861812592
$ws.Range("D2").Value = 0.03195955312829257
$ws.Range("E2").Value = 0.9690302269779759

$ws.Range("A3").Value = 6
$ws.Range("B3").Value = 0.3009774948852012
$ws.Range("C3").Value = 0.2567024827003479
$ws.Range("D3").Value = 0.0442750121848533
$ws.Range("E3").Value = 0.9576021530073574

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 0.106615139804501
$ws.Range("C4").Value = 0.1699984669685364
$ws.Range("D4").Value = 0.06338332716403534
$ws.Range("E4").Value = 0.9403946577448473

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 0.09297567629006591
$ws.Range("C5").Value = 0.1655720621347427
$ws.Range("D5").Value = 0.07259638584467683
$ws.Range("E5").Value = 0.9323171448246987

$ws.Range("A6").Value = 0
$ws.Range("B6").Value = 0.238463287110707
$ws.Range("C6").Value = 0.3150076270103455
$ws.Range("D6").Value = 0.07654433989963849
$ws.Range("E6").Value = 0.9288981075254418

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 0.2498295067060696
$ws.Range("C7").Value = 0.3265554904937744
$ws.Range("D7").Value = 0.07672598378770479
$ws.Range("E7").Value = 0.9287414022295643

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 0.9204364628324617
$ws.Range("C8").Value = 0.8226803541183472
$ws.Range("D8").Value = 0.09775610871411455
$ws.Range("E8").Value = 0.9109491553377701

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 0.1157081154807911
$ws.Range("C9").Value = 0.214569479227066
$ws.Range("D9").Value = 0.09886136374627495
$ws.Range("E9").Value = 0.9100329058715528
